$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: remove the last column (BA) entirely - shifts nothing else, just drops it
$ws.Columns.Item(53).Delete() | Out-Null

# Step 2: remove the last two rows (24, 23) entirely
$ws.Rows.Item(24).Delete() | Out-Null
$ws.Rows.Item(23).Delete() | Out-Null

# Step 3: update the recalculated forecast values that changed
$ws.Range("B1").Value = 39583
$ws.Range("C1").Value = 39765
$ws.Range("D1").Value = 39948
$ws.Range("E1").Value = 40130
$ws.Range("F1").Value = 40310
$ws.Range("G1").Value = 40494
$ws.Range("H1").Value = 40676
$ws.Range("I1").Value = 40862
$ws.Range("J1").Value = 41044
$ws.Range("K1").Value = 41228
$ws.Range("L1").Value = 41409
$ws.Range("M1").Value = 41592
$ws.Range("N1").Value = 41774
$ws.Range("O1").Value = 41957
$ws.Range("P1").Value = 42137
$ws.Range("Q1").Value = 42321
$ws.Range("R1").Value = 42503
$ws.Range("S1").Value = 42689
$ws.Range("T1").Value = 42867
$ws.Range("U1").Value = 43053
$ws.Range("V1").Value = 43145
$ws.Range("W1").Value = 43235
$ws.Range("X1").Value = 43326
$ws.Range("Y1").Value = 43418
$ws.Range("Z1").Value = 43510
$ws.Range("AA1").Value = 43600
$ws.Range("AB1").Value = 43691
$ws.Range("AC1").Value = 43783
$ws.Range("AD1").Value = 43875
$ws.Range("AE1").Value = 43966
$ws.Range("AF1").Value = 44068
$ws.Range("AG1").Value = 44159
$ws.Range("AH1").Value = 44251
$ws.Range("AI1").Value = 44341
$ws.Range("AJ1").Value = 44432
$ws.Range("AK1").Value = 44525
$ws.Range("AL1").Value = 44617
$ws.Range("AM1").Value = 44706
$ws.Range("AN1").Value = 44798
$ws.Range("AO1").Value = 44890
$ws.Range("AP1").Value = 44981
$ws.Range("AQ1").Value = 45071
$ws.Range("AR1").Value = 45163
$ws.Range("AS1").Value = 45254
$ws.Range("AT1").Value = 45345
$ws.Range("AU1").Value = 45436
$ws.Range("AV1").Value = 45534
$ws.Range("AW1").Value = 45618
$ws.Range("AX1").Value = 45713
$ws.Range("AY1").Value = 45800
$ws.Range("AZ1").Value = 45891
$ws.Range("B3").Value = 0.07975464681371225
$ws.Range("C3").Value = -4.700509864312973
$ws.Range("D3").Value = -3.017636378339217
$ws.Range("B4").Value = -0.1399818631928951
$ws.Range("C4").Value = -4.391509335919141
$ws.Range("D4").Value = -3.229247082222797
$ws.Range("E4").Value = -0.01655958389530365
$ws.Range("F4").Value = 0.003352386816724007
$ws.Range("D5").Value = -3.304481657602698
$ws.Range("E5").Value = -1.426203430357909
$ws.Range("F5").Value = -1.461031976610316
$ws.Range("G5").Value = 3.579142225970444
$ws.Range("H5").Value = 3.925837669383347
$ws.Range("F6").Value = -1.43036035661499
$ws.Range("G6").Value = 0.8577584548236317
$ws.Range("H6").Value = 2.641604203902781
$ws.Range("I6").Value = -0.289184878867832
$ws.Range("J6").Value = 0.2381541440396262
$ws.Range("H7").Value = 2.792069574291456
$ws.Range("I7").Value = 1.821104206634194
$ws.Range("J7").Value = 1.60268309892857
$ws.Range("K7").Value = 5.963492031746176
$ws.Range("L7").Value = 4.993892964711621
$ws.Range("J8").Value = 1.560945935618463
$ws.Range("K8").Value = 4.059266849997001
$ws.Range("L8").Value = 2.260118192030736
$ws.Range("M8").Value = 7.523777575896196
$ws.Range("N8").Value = 6.711795724673664
$ws.Range("L9").Value = 2.268053155954486
$ws.Range("M9").Value = 4.921660622329105
$ws.Range("N9").Value = 6.409878804372982
$ws.Range("O9").Value = 2.532215190177589
$ws.Range("P9").Value = 0.5121603413743347
$ws.Range("N10").Value = 6.652469936617145
$ws.Range("O10").Value = 5.215795625989261
$ws.Range("P10").Value = 3.290935868252554
$ws.Range("Q10").Value = 2.051185924063259
$ws.Range("R10").Value = 1.745565778643887
$ws.Range("P11").Value = 3.309637709230717
$ws.Range("Q11").Value = 2.840854095773526
$ws.Range("R11").Value = 0.7985845180024986
$ws.Range("S11").Value = 0.4575538530338541
$ws.Range("T11").Value = 2.687500891103922
$ws.Range("R12").Value = 0.6768288812109668
$ws.Range("S12").Value = 0.1992978909398646
$ws.Range("T12").Value = 1.922191950024699
$ws.Range("U12").Value = 2.600569166164624
$ws.Range("V12").Value = 3.339205815020496
$ws.Range("W12").Value = 3.654655474034474
$ws.Range("X12").Value = 4.036117574265741
$ws.Range("T13").Value = 1.852325089269979
$ws.Range("U13").Value = 2.217217717421827
$ws.Range("V13").Value = 2.646788941483735
$ws.Range("W13").Value = 3.068403604789749
$ws.Range("X13").Value = 3.628675245064317
$ws.Range("Y13").Value = 3.605726003451304
$ws.Range("Z13").Value = 3.813466308501412
$ws.Range("AA13").Value = 3.712036718632117
$ws.Range("AB13").Value = 3.551357200054261
$ws.Range("W14").Value = 3.010470130271137
$ws.Range("X14").Value = 3.259983363993291
$ws.Range("Y14").Value = 3.296463213734779
$ws.Range("Z14").Value = 3.884173085820986
$ws.Range("AA14").Value = 3.908921577463587
$ws.Range("AB14").Value = 3.524277826276134
$ws.Range("AC14").Value = 3.490656491795074
$ws.Range("AD14").Value = 2.732790977059629
$ws.Range("AE14").Value = 2.849400388885992
$ws.Range("AF14").Value = -0.985458715495402
$ws.Range("AA15").Value = 3.950347314142233
$ws.Range("AB15").Value = 3.648977694024791
$ws.Range("AC15").Value = 3.685313717535621
$ws.Range("AD15").Value = 3.412957258051663
$ws.Range("AE15").Value = 3.535456592693387
$ws.Range("AF15").Value = -3.099271113627677
$ws.Range("AG15").Value = -2.347097924577757
$ws.Range("AH15").Value = -5.285337128797329
$ws.Range("AI15").Value = -4.741003096464214
$ws.Range("AJ15").Value = -4.365687260408224
$ws.Range("AE16").Value = 3.508868582430846
$ws.Range("AF16").Value = 0.07095532091387913
$ws.Range("AG16").Value = -0.1565465363420615
$ws.Range("AH16").Value = -2.69436300797079
$ws.Range("AI16").Value = -2.156795995006056
$ws.Range("AJ16").Value = -0.9191921099315992
$ws.Range("AK16").Value = -0.1803381976702711
$ws.Range("AL16").Value = 1.871837441670499
$ws.Range("AM16").Value = 1.194925448553708
$ws.Range("AN16").Value = 1.27347919322387
$ws.Range("AH17").Value = -1.89286697317308
$ws.Range("AI17").Value = -1.937045172389718
$ws.Range("AJ17").Value = -1.824487515646256
$ws.Range("AK17").Value = -1.96738848373067
$ws.Range("AL17").Value = -0.8573220257725445
$ws.Range("AM17").Value = -2.068675356622807
$ws.Range("AN17").Value = -2.137023292796481
$ws.Range("AO17").Value = -1.152671696465724
$ws.Range("AP17").Value = -0.507642354784088
$ws.Range("AQ17").Value = -1.084365158506884
$ws.Range("AR17").Value = -1.339436245206127
$ws.Range("AL18").Value = -1.461464090310793
$ws.Range("AM18").Value = -2.197497829006645
$ws.Range("AN18").Value = -2.674188284733892
$ws.Range("AO18").Value = -0.5008759058252599
$ws.Range("AP18").Value = 0.07646803116447831
$ws.Range("AQ18").Value = -1.089896342664354
$ws.Range("AR18").Value = -1.277727682704721
$ws.Range("AS18").Value = -2.785556326028149
$ws.Range("AT18").Value = -2.87408779878463
$ws.Range("AU18").Value = -3.40787540386569
$ws.Range("AV18").Value = -3.451527003230626
$ws.Range("AP19").Value = -0.08815021972816695
$ws.Range("AQ19").Value = -0.8703448702657268
$ws.Range("AR19").Value = 0.07721345580697925
$ws.Range("AS19").Value = -0.6863115373258788
$ws.Range("AT19").Value = -1.155307395925487
$ws.Range("AU19").Value = -2.191935020614488
$ws.Range("AV19").Value = -2.998472503487815
$ws.Range("AW19").Value = -2.452009576682213
$ws.Range("AX19").Value = -2.164008261629446
$ws.Range("AY19").Value = -1.853660925652212
$ws.Range("AZ19").Value = -1.75044229618867
$ws.Range("AT20").Value = -1.187892669869473
$ws.Range("AU20").Value = -1.841569468248938
$ws.Range("AV20").Value = -2.923113274870115
$ws.Range("AW20").Value = -2.440246029655901
$ws.Range("AX20").Value = -2.337170009804157
$ws.Range("AY20").Value = -1.878672029998096
$ws.Range("AZ20").Value = -1.760724207457021
$ws.Range("AX21").Value = -2.433207997606113
$ws.Range("AY21").Value = -2.238303895464766
$ws.Range("AZ21").Value = -2.546230689156992

# Step 4: clear cells that are no longer part of the (now shorter) staircase pattern
$ws.Range("C5").ClearContents() | Out-Null
$ws.Range("E6").ClearContents() | Out-Null
$ws.Range("G7").ClearContents() | Out-Null
$ws.Range("I8").ClearContents() | Out-Null
$ws.Range("K9").ClearContents() | Out-Null
$ws.Range("M10").ClearContents() | Out-Null
$ws.Range("O11").ClearContents() | Out-Null
$ws.Range("Q12").ClearContents() | Out-Null
$ws.Range("R13").ClearContents() | Out-Null
$ws.Range("S13").ClearContents() | Out-Null
$ws.Range("T14").ClearContents() | Out-Null
$ws.Range("U14").ClearContents() | Out-Null
$ws.Range("V14").ClearContents() | Out-Null
$ws.Range("V15").ClearContents() | Out-Null
$ws.Range("W15").ClearContents() | Out-Null
$ws.Range("X15").ClearContents() | Out-Null
$ws.Range("Y15").ClearContents() | Out-Null
$ws.Range("Z15").ClearContents() | Out-Null
$ws.Range("Y16").ClearContents() | Out-Null
$ws.Range("Z16").ClearContents() | Out-Null
$ws.Range("AA16").ClearContents() | Out-Null
$ws.Range("AB16").ClearContents() | Out-Null
$ws.Range("AC16").ClearContents() | Out-Null
$ws.Range("AD16").ClearContents() | Out-Null
$ws.Range("AC17").ClearContents() | Out-Null
$ws.Range("AD17").ClearContents() | Out-Null
$ws.Range("AE17").ClearContents() | Out-Null
$ws.Range("AF17").ClearContents() | Out-Null
$ws.Range("AG17").ClearContents() | Out-Null
$ws.Range("AG18").ClearContents() | Out-Null
$ws.Range("AH18").ClearContents() | Out-Null
$ws.Range("AI18").ClearContents() | Out-Null
$ws.Range("AJ18").ClearContents() | Out-Null
$ws.Range("AK18").ClearContents() | Out-Null
$ws.Range("AK19").ClearContents() | Out-Null
$ws.Range("AL19").ClearContents() | Out-Null
$ws.Range("AM19").ClearContents() | Out-Null
$ws.Range("AN19").ClearContents() | Out-Null
$ws.Range("AO19").ClearContents() | Out-Null
$ws.Range("AO20").ClearContents() | Out-Null
$ws.Range("AP20").ClearContents() | Out-Null
$ws.Range("AQ20").ClearContents() | Out-Null
$ws.Range("AR20").ClearContents() | Out-Null
$ws.Range("AS20").ClearContents() | Out-Null
$ws.Range("AS21").ClearContents() | Out-Null
$ws.Range("AT21").ClearContents() | Out-Null
$ws.Range("AU21").ClearContents() | Out-Null
$ws.Range("AV21").ClearContents() | Out-Null
$ws.Range("AW21").ClearContents() | Out-Null
$ws.Range("AW22").ClearContents() | Out-Null
$ws.Range("AX22").ClearContents() | Out-Null
$ws.Range("AY22").ClearContents() | Out-Null
$ws.Range("AZ22").ClearContents() | Out-Null

Write-Host "edit complete"